# The deck's theme (ppt/theme/theme1.xml, used by the slide master) was the
# "Integral" / "Red Violet" palette.  It is being swapped for the default
# "Office Theme" palette (the one previously only used by the notes master).
#
# PowerPoint doesn't let you edit a theme's XML directly through the object
# model, but ThemeColorScheme.Colors(i).RGB is writable and maps 1:1 onto the
# <a:clrScheme> children (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in
# theme order, so we use it to push the "Office" RGB values onto the slide
# master's theme color scheme.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

# Office Theme color scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeTheme = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Colors($i).RGB = $officeTheme[$i - 1]
}
